$wb = $excel.ActiveWorkbook

# --- Rename the 4th sheet "Validações" -> "Dispositivos" ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "Dispositivos"

# --- Add a new sheet after all current sheets, named "Validaçoes" ---
$ws5 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws5.Name = "Validaçoes"

# Row 2
$ws5.Range("A2:A6").Merge()
$ws5.Range("A2:A6").HorizontalAlignment = -4108
$ws5.Range("A2:A6").VerticalAlignment = -4108
$ws5.Range("A2:A6").Orientation = 45
$ws5.Range("A2").Value = "validaçoes"

$ws5.Range("B2").Value = "Tipo"

$ws5.Range("C2:D2").Merge()
$ws5.Range("C2:D2").HorizontalAlignment = -4108
$ws5.Range("C2").Value = "W3c"

$ws5.Range("E2:H2").Merge()
$ws5.Range("E2:H2").HorizontalAlignment = -4108
$ws5.Range("E2:H2").Font.Underline = $true
$ws5.Range("E2").Value = "Google"

$ws5.Range("I2:L2").Merge()
$ws5.Range("I2:L2").HorizontalAlignment = -4108
$ws5.Range("I2").Value = "Acessibilidade"

# Row 3
$ws5.Range("B3").Value = "Status"
$ws5.Range("C3").Value = "Html"
$ws5.Range("D3").Value = "Css"

$ws5.Range("E3:F3").Merge()
$ws5.Range("E3:F3").HorizontalAlignment = -4131
$ws5.Range("E3:F3").VerticalAlignment = -4160
$ws5.Range("E3").Value = "segurança"

$ws5.Range("G3:H3").Merge()
$ws5.Range("G3:H3").HorizontalAlignment = -4131
$ws5.Range("G3").Value = "Velocidade"

# Row 4
$ws5.Range("B4").Value = "Data"

# Row 6 - empty underline-styled cell
$ws5.Range("L6").Font.Underline = $true

# Selection / active cell on the new sheet
$ws5.Range("L6").Select()
